$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''243.35'
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Value = '''23.78'
$ws.Range("D3").Style = "Normal"
$ws.Range("D4").Value = '''5.393'
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").Value = '''0.05884'
$ws.Range("D5").Style = "Normal"
$ws.Range("B6").Value = 'KuCoinToken'
$ws.Range("C6").Value = 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'
$ws.Range("D6").Value = '''6.504'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '5KuCoinTokenKCS'
$ws.Range("B7").Value = 'GateToken'
$ws.Range("C7").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D7").Value = '''3.383'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '6GateTokenGT'
$ws.Range("D8").Value = '''0.8114'
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Value = '''0.9271'
$ws.Range("D9").Style = "Normal"
$ws.Range("B10").Value = 'WazirX'
$ws.Range("C10").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D10").Value = '''0.1416'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '9WazirXWRX'
$ws.Range("B11").Value = 'MandalaExchangeToken'
$ws.Range("C11").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D11").Value = '''0.07389'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '10MandalaExchangeTokenMDX'
$ws.Range("B12").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C12").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D12").Value = '''0.03050'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '11LiechtensteinCryptoassetsExchangeLCX'
$ws.Range("B13").Value = 'BitrueCoin'
$ws.Range("C13").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D13").Value = '''0.03056'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '12BitrueCoinBTR'
$ws.Range("B14").Value = 'BitMartToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D14").Value = '''0.09346'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '13BitMartTokenBMX'
$ws.Range("B15").Value = 'MCDex'
$ws.Range("C15").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range("D15").Value = '''3.860'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '14MCDexMCB'
$ws.Range("B16").Value = 'BitForexToken'
$ws.Range("C16").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D16").Value = '''0.001554'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '15BitForexTokenBF'
$ws.Range("B17").Value = 'CoinExToken'
$ws.Range("C17").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range("D17").Value = '''0.04699'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '16CoinExTokenCET'
$ws.Range("B18").Value = 'One'
$ws.Range("C18").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D18").Value = '''0.0005986'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '17OneONE'
$ws.Range("D19").Value = '''0.005867'
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Value = '''0.001244'
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Value = '''0.004721'
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Value = '''0.00008805'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '21NitroExNTXBestin24h'
$ws.Range("D23").Value = '''3.556'
$ws.Range("D23").Style = "Normal"
$ws.Range("D25").Value = '''0.3227'
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Value = '''0.1330'
$ws.Range("D26").Style = "Normal"
$ws.Range("D40").Value = '''0.03870'
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Value = '''0.003139'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '40KickTokenKICKWorstin24h'
$ws.Range("B42").Value = 'BKEXToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range("D42").Value = '''0.1068'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '41BKEXTokenBKK'
$ws.Range("B43").Value = 'CEJI'
$ws.Range("C43").Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
$ws.Range("D43").Value = '''0.002802'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '42CEJICEJI'
$ws.Range("D44").Value = '''0.008570'
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Value = '''0.00005214'
$ws.Range("D45").Style = "Normal"
$ws.Range("D47").Value = '''0.6717'
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Value = '''0.001946'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '47BOLOBOLO'
